$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# NB: Shape.Left/.Width are expressed in points (1 pt = 12700 EMU) and are
# stored by the host as single-precision floats, then converted back to EMU
# by truncation. The literals below are the exact point values whose
# float32 representation reproduces the target EMU offsets/extents from the
# diff (5086173/2555684, 8158785/2555684, 11750289/2555684) once saved back
# to OOXML.

# --- Shape "Rectangle 10" (id 11) : "mouse hover" -> "3d mouse hover" ---
$sh1 = $s.Shapes.Item(1)
$sh1.Left = 400.486083984375
$sh1.Width = 201.23497009277344
$tr1 = $sh1.TextFrame2.TextRange
$tr1.Text = "3d mouse hover"
$tr1.Characters(1, 9).Text = "3d mouse "

# --- Shape "Rectangle 10" (id 13) : "mouse touch click" -> "3d mouse touch click" ---
$sh2 = $s.Shapes.Item(2)
$sh2.Left = 642.424072265625
$sh2.Width = 201.23497009277344
$tr2 = $sh2.TextFrame2.TextRange
$tr2.Text = "3d mouse touch click"
$tr2.Characters(1, 9).Text = "3d mouse "

# --- Shape "Rectangle 10" (id 17) : "mouse touch hold" -> "3d mouse touch hold" ---
$sh3 = $s.Shapes.Item(3)
$sh3.Left = 925.2196655273438
$sh3.Width = 201.23497009277344
$tr3 = $sh3.TextFrame2.TextRange
$tr3.Text = "3d mouse touch hold"
$tr3.Characters(1, 9).Text = "3d mouse "
